# Rename the embedded logo pictures (Word exposes this as InlineShape.Name,
# which maps to the picture's wp:docPr/@name OOXML attribute):
#   - First-page header BTec logo:      image1.jpg -> image2.jpg
#   - Default footer Pearson logo:      image2.png -> image1.png
#   - First-page footer Pearson logo:   image2.png -> image1.png

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# First-page header (wdHeaderFooterFirstPage) holds the BTec_Logo-Orange picture.
$firstHeader = $sec.Headers.Item(2)
if ($firstHeader.Exists -and $firstHeader.Range.InlineShapes.Count -gt 0) {
    $btecShape = $firstHeader.Range.InlineShapes.Item(1)
    if ($btecShape.AlternativeText -eq "BTec_Logo-Orange") {
        $btecShape.Name = "image2.jpg"
    }
}

# Default + first-page footers both hold the Pearson logo picture.
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        $pearsonShape = $ftr.Range.InlineShapes.Item(1)
        if ($pearsonShape.AlternativeText -like "*PearsonLogo.png") {
            $pearsonShape.Name = "image1.png"
        }
    }
}
